# Auto-generated edit script: update crypto price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) updates that are safely non-numeric text already ---
$ws.Range("D2").Value = "60.628.99"
$ws.Range("D3").Value = "2.402.87"
$ws.Range("D9").Value = "2.384.28"
$ws.Range("D15").Value = "2.843.33"
$ws.Range("D17").Value = "60.733.45"
$ws.Range("D18").Value = "2.405.41"
$ws.Range("D29").Value = "2.530.59"
$ws.Range("D45").Value = "0.0₆0289"

# --- Price (D) updates that look like plain numbers: force text type ---
# so Excel keeps them as inline/shared strings instead of converting to numeric,
# matching the original inlineStr cell type. We temporarily mark the cell as
# Text-formatted, assign the value, then restore the default "Normal" style so
# no visible formatting change remains on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "548.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "152.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "141.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.580"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0894"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) (E) updates: always text (percentages with padding spaces) ---
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -3.12%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("E19").Value = "  +10.44%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("E23").Value = "  -9.81%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -6.68%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -9.99%  "
$ws.Range("E28").Value = "  -5.95%  "
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  -5.54%  "
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("E39").Value = "  -5.27%  "
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("E45").Value = "  +4.78%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("E51").Value = "  -0.16%  "
